$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Logistic Regression row (row 2) with refreshed metrics
$ws.Range("B2").Value = 0.7691460055096419
$ws.Range("C2").Value = 0.7724018023761483
$ws.Range("D2").Value = 0.7691460055096419
$ws.Range("E2").Value = 0.7608422551532319

# Replace row 3 (was "Lasso") with the "LightGBM" row that used to be row 7
$ws.Range("A3").Value = "LightGBM"
$ws.Range("B3").Value = 0.796969696969697
$ws.Range("C3").Value = 0.7962065687395939
$ws.Range("D3").Value = 0.796969696969697
$ws.Range("E3").Value = 0.7937582430116912

# Remove the now-unused model rows (Support Vector Classifier, CART, Random Forest,
# old LightGBM, XGBoost) so only the two relevant models remain
$ws.Range("A4:E8").Delete()
